$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-07-15 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-07-16 Tuesday", 2) | Out-Null
$d.Content.Find.Execute("641÷3=213, 2", $true, $false, $false, $false, $false, $true, 1, $false, "238÷9=26, 4", 2) | Out-Null
$d.Content.Find.Execute("464÷4=116, 0", $true, $false, $false, $false, $false, $true, 1, $false, "358÷9=39, 7", 2) | Out-Null
$d.Content.Find.Execute("723÷6=120, 3", $true, $false, $false, $false, $false, $true, 1, $false, "888÷2=444, 0", 2) | Out-Null
$d.Content.Find.Execute("285÷6=47, 3", $true, $false, $false, $false, $false, $true, 1, $false, "898÷2=449, 0", 2) | Out-Null
$d.Content.Find.Execute("612÷5=122, 2", $true, $false, $false, $false, $false, $true, 1, $false, "888÷5=177, 3", 2) | Out-Null
$d.Content.Find.Execute("417÷6=69, 3", $true, $false, $false, $false, $false, $true, 1, $false, "389÷9=43, 2", 2) | Out-Null
$d.Content.Find.Execute("913÷9=101, 4", $true, $false, $false, $false, $false, $true, 1, $false, "389÷7=55, 4", 2) | Out-Null
$d.Content.Find.Execute("814÷3=271, 1", $true, $false, $false, $false, $false, $true, 1, $false, "938÷7=134, 0", 2) | Out-Null
$d.Content.Find.Execute("143÷7=20, 3", $true, $false, $false, $false, $false, $true, 1, $false, "169÷8=21, 1", 2) | Out-Null
$d.Content.Find.Execute("810÷6=135, 0", $true, $false, $false, $false, $false, $true, 1, $false, "119÷6=19, 5", 2) | Out-Null
$d.Content.Find.Execute("466÷4=116, 2", $true, $false, $false, $false, $false, $true, 1, $false, "277÷4=69, 1", 2) | Out-Null
$d.Content.Find.Execute("931÷2=465, 1", $true, $false, $false, $false, $false, $true, 1, $false, "157÷3=52, 1", 2) | Out-Null
$d.Content.Find.Execute("178÷8=22, 2", $true, $false, $false, $false, $false, $true, 1, $false, "649÷6=108, 1", 2) | Out-Null
$d.Content.Find.Execute("967÷4=241, 3", $true, $false, $false, $false, $false, $true, 1, $false, "256÷6=42, 4", 2) | Out-Null
$d.Content.Find.Execute("212÷9=23, 5", $true, $false, $false, $false, $false, $true, 1, $false, "269÷3=89, 2", 2) | Out-Null
$d.Content.Find.Execute("600÷2=300, 0", $true, $false, $false, $false, $false, $true, 1, $false, "376÷7=53, 5", 2) | Out-Null
$d.Content.Find.Execute("509÷2=254, 1", $true, $false, $false, $false, $false, $true, 1, $false, "510÷8=63, 6", 2) | Out-Null
$d.Content.Find.Execute("900÷2=450, 0", $true, $false, $false, $false, $false, $true, 1, $false, "100÷9=11, 1", 2) | Out-Null
$d.Content.Find.Execute("696÷7=99, 3", $true, $false, $false, $false, $false, $true, 1, $false, "250÷4=62, 2", 2) | Out-Null
$d.Content.Find.Execute("225÷3=75, 0", $true, $false, $false, $false, $false, $true, 1, $false, "335÷2=167, 1", 2) | Out-Null
$d.Content.Find.Execute("954÷6=159, 0", $true, $false, $false, $false, $false, $true, 1, $false, "496÷6=82, 4", 2) | Out-Null
$d.Content.Find.Execute("676÷5=135, 1", $true, $false, $false, $false, $false, $true, 1, $false, "140÷4=35, 0", 2) | Out-Null
$d.Content.Find.Execute("561÷2=280, 1", $true, $false, $false, $false, $false, $true, 1, $false, "881÷3=293, 2", 2) | Out-Null
$d.Content.Find.Execute("496÷4=124, 0", $true, $false, $false, $false, $false, $true, 1, $false, "479÷9=53, 2", 2) | Out-Null
$d.Content.Find.Execute("834÷2=417, 0", $true, $false, $false, $false, $false, $true, 1, $false, "103÷6=17, 1", 2) | Out-Null
